$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1): "<col>_old" -> "<col>_FV2410" (A..J)
#        and "<col>_new" -> "<col>_FV2504" (L..U). K1 ("diff") is unchanged.
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value2 = $fv2410Headers[$i]
}
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value2 = $fv2504Headers[$i]
}

# --- 2. Turn the data range into a native Excel Table (ListObject) named "Table1" ---
$rng = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $rng, [System.Type]::Missing, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (split/frozen pane below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

"done"
